$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds a small "TestCase" table with repeating storeText
# target/value column-pairs. This change inserts one more storeText
# pair (a new column) right before the existing "wait" column, so a
# second storeText step (id=div2 / value3) can be recorded.

# 1) Insert a new column before E; this shifts the old E ("wait") to F.
$ws.Columns("E:E").Insert()

# 2) Populate the new column's header (row 1) with the same "storeText"
#    label used by the other storeText column (D), and give it the same
#    column width as D.
$ws.Range("E1").Value = $ws.Range("D1").Value2
$ws.Columns("E:E").ColumnWidth = 23.65

# 3) Row 2 (target column) stays blank for the new storeText pair, just
#    like D2/E2 for the first pair -- nothing else to set there.

# 4) Row 3 holds the value for the new storeText step.
$ws.Range("E3").Value = '{"target":"id=div2","value":"value3"}'
$ws.Rows("3:3").RowHeight = 37.5

# 5) Drop the redundant direct format on C2 (it only set a no-op
#    "apply fill" flag on top of the default format) by pasting in a
#    clean default format from an untouched cell.
$ws.Range("Z1").ClearFormats()
$ws.Range("Z1").Copy()
$ws.Range("C2").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("Z1").Clear()

# 6) Match the saved selection/cursor position.
$ws.Range("E5").Select()
